# Completed excel sheet for investment_account
# Fill in developer name, preconditions, method inputs and expected
# results for the InvestmentAccount unit-test plan (Table1, rows 7-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Sahil Choudhary"

# Preconditions (column E) - "None" everywhere, first used on row 7
$ws.Range("E7").Value = "None"

# Expected Result (column G) for rows 7-10
$ws.Range("G7").Value = "Attributes are set"
$ws.Range("G8").Value = "management_fee set to 2.55"
$ws.Range("G9").Value = "service_charge set to 0.50"
$ws.Range("G10").Value = "service_charge set to 2.50"

# Method Inputs (column F) for rows 7, 9, 8, 10
$ws.Range("F7").Value = "account_number = 350, client_number = 350, balance = 350, date_created = (2024, 3, 25), management_fee = 2"
$ws.Range("F9").Value = "account_number = 350, client_number = 350, balance = 350, date_created = (2014, 3, 25), management_fee = 2"
$ws.Range("F8").Value = "account_number = 350, client_number = 350, balance = 350, date_created = (2024, 3, 25), management_fee = 'two'"
$ws.Range("F10").Value = "account_number = 350, client_number = 350, balance = 350, date_created = InvestmentAccount.TEN_YEARS_AGO, management_fee = 2"

# Expected Result (column G) for rows 12, 13
$ws.Range("G12").Value = "Account Number: 350 Balance: $350.00 Date Created: 2014-03-25 Management Fee: Waived Account Type: Investment"
$ws.Range("G13").Value = "Account Number: 350 Balance: $350.00 Date Created: 2024-03-25 Management Fee: $2.00 Account Type: Investment"

# Remaining cells that reuse already-existing shared strings
$ws.Range("E8").Value = "None"
$ws.Range("E9").Value = "None"
$ws.Range("E10").Value = "None"
$ws.Range("E11").Value = "None"
$ws.Range("E12").Value = "None"
$ws.Range("E13").Value = "None"

$ws.Range("F11").Value = "account_number = 350, client_number = 350, balance = 350, date_created = (2024, 3, 25), management_fee = 2"
$ws.Range("F12").Value = "account_number = 350, client_number = 350, balance = 350, date_created = (2014, 3, 25), management_fee = 2"
$ws.Range("F13").Value = "account_number = 350, client_number = 350, balance = 350, date_created = (2024, 3, 25), management_fee = 2"

$ws.Range("G11").Value = "service_charge set to 2.50"

# Final selection left on the sheet (matches author's last-saved cursor position)
$ws.Range("G13").Select()
